$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1917.5834
$ws.Range("I33").Value = 1281.1
$ws.Range("J33").Value = 5100
$ws.Range("K33").Value = 1281.1
$ws.Range("L33").Value = 5100
$ws.Range("M33").Value = -1052.1
$ws.Range("N33").Value = -5558
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H113").Value = 2338.25
$ws.Range("I113").Value = 2140
$ws.Range("J113").Value = 2668.6667
$ws.Range("K113").Value = 2140
$ws.Range("L113").Value = 2668.6667
$ws.Range("M113").Value = 1114
$ws.Range("N113").Value = -9176.6667
$ws.Range("H137").Value = 3309.7693
$ws.Range("I137").Value = 2359.7334
$ws.Range("J137").Value = 6476.5557
$ws.Range("K137").Value = 7079.2002
$ws.Range("L137").Value = 19429.6671
$ws.Range("M137").Value = -4529.2002
$ws.Range("N137").Value = -24529.6671
$ws.Range("H138").Value = 2050.84
$ws.Range("J138").Value = 2237.3296
$ws.Range("L138").Value = 6711.9888
$ws.Range("N138").Value = -16991.9888

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 73600.78999999999
$ws.Range("I122").Value = 112745.664
$ws.Range("J122").Value = 3140
$ws.Range("K122").Value = 338236.992
$ws.Range("L122").Value = 9420
$ws.Range("M122").Value = -335786.992
$ws.Range("N122").Value = -14320
$ws.Range("H132").Value = 4317.5527
$ws.Range("I132").Value = 3650.9443
$ws.Range("J132").Value = 4917.5
$ws.Range("K132").Value = 10952.8329
$ws.Range("L132").Value = 14752.5
$ws.Range("M132").Value = -8422.832900000001
$ws.Range("N132").Value = -19812.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 142860060
$ws.Range("I86").Value = 142860060
$ws.Range("K86").Value = 142860060
$ws.Range("M86").Value = -142858937
$ws.Range("H89").Value = 142860060
$ws.Range("I89").Value = 142860060
$ws.Range("K89").Value = 714300300
$ws.Range("M89").Value = -714294684
$ws.Range("H134").Value = 2503.111
$ws.Range("I134").Value = 2416.9092
$ws.Range("J134").Value = 2638.5715
$ws.Range("K134").Value = 7250.7276
$ws.Range("L134").Value = 7915.7145
$ws.Range("M134").Value = -4715.7276
$ws.Range("N134").Value = -12985.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 11799
$ws.Range("J51").Value = 11799
$ws.Range("L51").Value = 11799
$ws.Range("N51").Value = -13271
$ws.Range("H59").Value = 16799.2
$ws.Range("J59").Value = 16799.2
$ws.Range("L59").Value = 16799.2
$ws.Range("N59").Value = -19089.2
$ws.Range("H60").Value = 9858.429
$ws.Range("J60").Value = 9858.429
$ws.Range("L60").Value = 9858.429
$ws.Range("N60").Value = -10880.429
$ws.Range("H61").Value = 11799
$ws.Range("J61").Value = 11799
$ws.Range("L61").Value = 11799
$ws.Range("N61").Value = -12495
$ws.Range("H68").Value = 23518
$ws.Range("J68").Value = 23518
$ws.Range("L68").Value = 23518
$ws.Range("N68").Value = -25016
$ws.Range("H71").Value = 23518
$ws.Range("J71").Value = 23518
$ws.Range("L71").Value = 70554
$ws.Range("N71").Value = -78042
$ws.Range("H74").Value = 26258.455
$ws.Range("J74").Value = 26258.455
$ws.Range("L74").Value = 26258.455
$ws.Range("N74").Value = -28006.455
$ws.Range("H77").Value = 26258.455
$ws.Range("J77").Value = 26258.455
$ws.Range("L77").Value = 78775.36500000001
$ws.Range("N77").Value = -87511.36500000001
$ws.Range("H107").Value = 3290303.2
$ws.Range("I107").Value = 5682573.5
$ws.Range("J107").Value = 931.25
$ws.Range("K107").Value = 5682573.5
$ws.Range("L107").Value = 931.25
$ws.Range("M107").Value = -5680653.5
$ws.Range("N107").Value = -4771.25
$ws.Range("H120").Value = 26666.666
$ws.Range("J120").Value = 26666.666
$ws.Range("L120").Value = 26666.666
$ws.Range("N120").Value = -33924.666
$ws.Range("H122").Value = 1863.1364
$ws.Range("I122").Value = 1250.25
$ws.Range("J122").Value = 1999.3334
$ws.Range("K122").Value = 3750.75
$ws.Range("L122").Value = 5998.0002
$ws.Range("M122").Value = -1300.75
$ws.Range("N122").Value = -10898.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1541.3334
$ws.Range("I38").Value = 244
$ws.Range("J38").Value = 2190
$ws.Range("K38").Value = 732
$ws.Range("L38").Value = 6570
$ws.Range("M38").Value = -385
$ws.Range("N38").Value = -7264
$ws.Range("H39").Value = 1893.8
$ws.Range("J39").Value = 1893.8
$ws.Range("L39").Value = 5681.4
$ws.Range("N39").Value = -6269.4
$ws.Range("H107").Value = 576.44446
$ws.Range("I107").Value = 381.33334
$ws.Range("J107").Value = 674
$ws.Range("K107").Value = 1144.00002
$ws.Range("L107").Value = 2022
$ws.Range("M107").Value = 775.9999800000001
$ws.Range("N107").Value = -5862
$ws.Range("H110").Value = 12212.071
$ws.Range("I110").Value = 4513.5
$ws.Range("J110").Value = 12804.27
$ws.Range("K110").Value = 13540.5
$ws.Range("L110").Value = 38412.81
$ws.Range("M110").Value = -9450.5
$ws.Range("N110").Value = -46592.81
$ws.Range("H113").Value = 1227.1177
$ws.Range("J113").Value = 1248.8125
$ws.Range("L113").Value = 3746.4375
$ws.Range("N113").Value = -8086.4375
$ws.Range("H122").Value = 14828.286
$ws.Range("I122").Value = 360
$ws.Range("K122").Value = 3240
$ws.Range("M122").Value = -790
$ws.Range("H131").Value = 1011.8823
$ws.Range("J131").Value = 1033.7959
$ws.Range("L131").Value = 3101.3877
$ws.Range("N131").Value = -13181.3877
$ws.Range("H134").Value = 4273.6665
$ws.Range("I134").Value = 3094.0908
$ws.Range("J134").Value = 7517.5
$ws.Range("K134").Value = 9282.2724
$ws.Range("L134").Value = 22552.5
$ws.Range("M134").Value = -4212.2724
$ws.Range("N134").Value = -32692.5
$ws.Range("H137").Value = 6180620
$ws.Range("J137").Value = 4227.4287
$ws.Range("L137").Value = 12682.2861
$ws.Range("N137").Value = -22882.2861
$ws.Range("H139").Value = 2985.5
$ws.Range("I139").Value = 1728.2222
$ws.Range("J139").Value = 3739.8667
$ws.Range("K139").Value = 5184.6666
$ws.Range("L139").Value = 11219.6001
$ws.Range("M139").Value = -44.66659999999956
$ws.Range("N139").Value = -21499.6001
$ws.Range("H141").Value = 5955.2
$ws.Range("I141").Value = 1446.25
$ws.Range("J141").Value = 8077.0586
$ws.Range("K141").Value = 4338.75
$ws.Range("L141").Value = 24231.1758
$ws.Range("M141").Value = 841.25
$ws.Range("N141").Value = -34591.1758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5784.9395
$ws.Range("I70").Value = 5853.095
$ws.Range("J70").Value = 5665.6665
$ws.Range("K70").Value = 5853.095
$ws.Range("L70").Value = 5665.6665
$ws.Range("M70").Value = -5583.095
$ws.Range("N70").Value = -6205.6665
$ws.Range("H73").Value = 5784.9395
$ws.Range("I73").Value = 5853.095
$ws.Range("J73").Value = 5665.6665
$ws.Range("K73").Value = 5853.095
$ws.Range("L73").Value = 5665.6665
$ws.Range("M73").Value = -4917.095
$ws.Range("N73").Value = -7537.6665
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H107").Value = 241.71428
$ws.Range("I107").Value = 200.4
$ws.Range("J107").Value = 345
$ws.Range("K107").Value = 200.4
$ws.Range("L107").Value = 345
$ws.Range("M107").Value = 1719.6
$ws.Range("N107").Value = -4185
$ws.Range("H113").Value = 251049.75
$ws.Range("I113").Value = 1000000
$ws.Range("J113").Value = 1399.6666
$ws.Range("K113").Value = 1000000
$ws.Range("L113").Value = 1399.6666
$ws.Range("M113").Value = -997830
$ws.Range("N113").Value = -5739.6666
$ws.Range("H126").Value = 1996.625
$ws.Range("I126").Value = 1994.3334
$ws.Range("J126").Value = 2003.5
$ws.Range("K126").Value = 5983.0002
$ws.Range("L126").Value = 6010.5
$ws.Range("M126").Value = -3513.0002
$ws.Range("N126").Value = -10950.5
$ws.Range("H132").Value = 3626
$ws.Range("I132").Value = 3216.2856
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 9648.856800000001
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -7118.856800000001
$ws.Range("N132").Value = -17658.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4942.125
$ws.Range("I61").Value = 4853.857
$ws.Range("J61").Value = 5560
$ws.Range("K61").Value = 4853.857
$ws.Range("L61").Value = 5560
$ws.Range("M61").Value = -4651.857
$ws.Range("N61").Value = -5964
$ws.Range("H95").Value = 45086
$ws.Range("J95").Value = 45086
$ws.Range("L95").Value = 45086
$ws.Range("N95").Value = -50578
$ws.Range("H113").Value = 4942.125
$ws.Range("I113").Value = 4853.857
$ws.Range("J113").Value = 5560
$ws.Range("K113").Value = 4853.857
$ws.Range("L113").Value = 5560
$ws.Range("M113").Value = -2683.857
$ws.Range("N113").Value = -9900
$ws.Range("H122").Value = 3198.35
$ws.Range("I122").Value = 3088.3572
$ws.Range("J122").Value = 3455
$ws.Range("K122").Value = 9265.071599999999
$ws.Range("L122").Value = 10365
$ws.Range("M122").Value = -6815.071599999999
$ws.Range("N122").Value = -15265
$ws.Range("H136").Value = 9261502
$ws.Range("I136").Value = 1821.875
$ws.Range("J136").Value = 16669245
$ws.Range("K136").Value = 5465.625
$ws.Range("L136").Value = 50007735
$ws.Range("M136").Value = -2915.625
$ws.Range("N136").Value = -50012835
$ws.Range("H140").Value = 48809.332
$ws.Range("J140").Value = 48809.332
$ws.Range("L140").Value = 48809.332
$ws.Range("N140").Value = -59169.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 67125
$ws.Range("J97").Value = 67125
$ws.Range("L97").Value = 67125
$ws.Range("N97").Value = -69107
$ws.Range("H122").Value = 2904
$ws.Range("I122").Value = 1992.8
$ws.Range("J122").Value = 3663.3333
$ws.Range("K122").Value = 5978.4
$ws.Range("L122").Value = 10989.9999
$ws.Range("M122").Value = -3528.4
$ws.Range("N122").Value = -15889.9999
$ws.Range("H132").Value = 4506942
$ws.Range("I132").Value = 2752.7646
$ws.Range("J132").Value = 8335503
$ws.Range("K132").Value = 8258.293799999999
$ws.Range("L132").Value = 25006509
$ws.Range("M132").Value = -5728.293799999999
$ws.Range("N132").Value = -25011569
$ws.Range("H136").Value = 3097.9333
$ws.Range("I136").Value = 3206.2778
$ws.Range("J136").Value = 2935.4167
$ws.Range("K136").Value = 9618.8334
$ws.Range("L136").Value = 8806.250100000001
$ws.Range("M136").Value = -7068.8334
$ws.Range("N136").Value = -13906.2501
